# "Generate Report for Handoff"
# The localization report got regenerated: the d90b8816 e2e file moved from
# "Ready for handoff" (row 7) to "In Translation" (row 6) with a fresh
# handoff timestamp, while b7277af2 (previously row 6) drops to row 7 and
# keeps its "Ready for handoff" status. This swaps the two rows' content on
# every sheet (Overview, zh-cn, de-de) and updates the hyperlink display text
# accordingly (the hyperlink targets/rIds themselves are left untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "d90b8816-4b3a-47be-9bdf-38b007fe90ff.md"
$ws.Range("B6").Value = "e2e\d90b8816-4b3a-47be-9bdf-38b007fe90ff.md"
$ws.Range("C6").Value = ".md"
$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"
$ws.Range("G6").Value = "2016-11-07 06:58:48"

$ws.Range("A7").Value = "b7277af2-cbe8-4e42-9c08-0f7360d1714a.md"
$ws.Range("B7").Value = "e2e\b7277af2-cbe8-4e42-9c08-0f7360d1714a.md"
$ws.Range("C7").Value = ".md"
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-11-07 06:51:08"

# Rebuild the hyperlinks so the display text follows the swapped rows while
# keeping the same targets (deleting any single hyperlink clears the whole
# collection in this engine, so all of them are re-added together).
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/0c47f287-8800-41bb-bee0-29a00708bd1b.md", "", "", "e2e\0c47f287-8800-41bb-bee0-29a00708bd1b.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/209c48c6-fc24-415d-a9e7-73c5a75a6102.md", "", "", "e2e\209c48c6-fc24-415d-a9e7-73c5a75a6102.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16fdc2acbe502e300bb1fd34d7e148d58d66d8d9/e2e/70d2e4f3-9044-424f-b315-358cf7ddc94d.md", "", "", "e2e\70d2e4f3-9044-424f-b315-358cf7ddc94d.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7100a5fa181df88894f3deac028304d24a80d08b/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md", "", "", "e2e\a14b8ca5-f559-4148-9701-350adb07cd9e.md")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0278267dba7863fb2271ceb76abd24e533f8da94/e2e/b7277af2-cbe8-4e42-9c08-0f7360d1714a.md", "", "", "e2e\d90b8816-4b3a-47be-9bdf-38b007fe90ff.md")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5b2ff6e78250a65886f6d39eb837ad365ccb970/e2e/d90b8816-4b3a-47be-9bdf-38b007fe90ff.md", "", "", "e2e\b7277af2-cbe8-4e42-9c08-0f7360d1714a.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A6").Value = "d90b8816-4b3a-47be-9bdf-38b007fe90ff.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("G6").Value = "d90b8816-4b3a-47be-9bdf-38b007fe90ff.469eb875979ee048aaa47db77e61783c3208b67c.zh-cn.xlf"
$ws.Range("H6").Value = "2016-11-07 06:58:34"

$ws.Range("A7").Value = "b7277af2-cbe8-4e42-9c08-0f7360d1714a.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "b7277af2-cbe8-4e42-9c08-0f7360d1714a.48151d32738d36ccf46c134c39835abf022fbe6c.zh-cn.xlf"
$ws.Range("H7").Value = "2016-11-07 06:50:54"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/0c47f287-8800-41bb-bee0-29a00708bd1b.md", "", "", "0c47f287-8800-41bb-bee0-29a00708bd1b.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/209c48c6-fc24-415d-a9e7-73c5a75a6102.md", "", "", "209c48c6-fc24-415d-a9e7-73c5a75a6102.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16fdc2acbe502e300bb1fd34d7e148d58d66d8d9/e2e/70d2e4f3-9044-424f-b315-358cf7ddc94d.md", "", "", "70d2e4f3-9044-424f-b315-358cf7ddc94d.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7100a5fa181df88894f3deac028304d24a80d08b/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md", "", "", "a14b8ca5-f559-4148-9701-350adb07cd9e.md")
$ws.Hyperlinks.Add($ws.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/657eef637aaca1f94ddfc8c24d2604d3f9e5b5c8/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md", "", "", "a14b8ca5-f559-4148-9701-350adb07cd9e.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0278267dba7863fb2271ceb76abd24e533f8da94/e2e/b7277af2-cbe8-4e42-9c08-0f7360d1714a.md", "", "", "d90b8816-4b3a-47be-9bdf-38b007fe90ff.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5b2ff6e78250a65886f6d39eb837ad365ccb970/e2e/d90b8816-4b3a-47be-9bdf-38b007fe90ff.md", "", "", "b7277af2-cbe8-4e42-9c08-0f7360d1714a.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A6").Value = "d90b8816-4b3a-47be-9bdf-38b007fe90ff.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("G6").Value = "d90b8816-4b3a-47be-9bdf-38b007fe90ff.469eb875979ee048aaa47db77e61783c3208b67c.de-de.xlf"
$ws.Range("H6").Value = "2016-11-07 06:58:48"

$ws.Range("A7").Value = "b7277af2-cbe8-4e42-9c08-0f7360d1714a.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "b7277af2-cbe8-4e42-9c08-0f7360d1714a.48151d32738d36ccf46c134c39835abf022fbe6c.de-de.xlf"
$ws.Range("H7").Value = "2016-11-07 06:51:08"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/0c47f287-8800-41bb-bee0-29a00708bd1b.md", "", "", "0c47f287-8800-41bb-bee0-29a00708bd1b.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5695de6e623af8e9c7658d63e9485da1d4679/e2e/209c48c6-fc24-415d-a9e7-73c5a75a6102.md", "", "", "209c48c6-fc24-415d-a9e7-73c5a75a6102.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16fdc2acbe502e300bb1fd34d7e148d58d66d8d9/e2e/70d2e4f3-9044-424f-b315-358cf7ddc94d.md", "", "", "70d2e4f3-9044-424f-b315-358cf7ddc94d.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7100a5fa181df88894f3deac028304d24a80d08b/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md", "", "", "a14b8ca5-f559-4148-9701-350adb07cd9e.md")
$ws.Hyperlinks.Add($ws.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8081c555dd32f546b23969e8d22e3ed498a82946/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md", "", "", "a14b8ca5-f559-4148-9701-350adb07cd9e.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0278267dba7863fb2271ceb76abd24e533f8da94/e2e/b7277af2-cbe8-4e42-9c08-0f7360d1714a.md", "", "", "d90b8816-4b3a-47be-9bdf-38b007fe90ff.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5b2ff6e78250a65886f6d39eb837ad365ccb970/e2e/d90b8816-4b3a-47be-9bdf-38b007fe90ff.md", "", "", "b7277af2-cbe8-4e42-9c08-0f7360d1714a.md")
